$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M3: 5 -> 6
$ws.Range("M3").Value = "6"

# Row 4 and Row 5 swap their content (author-id corrected order),
# with row4's cited_by_count (M4) updated from 0 to 1.

# New row 4 (previously row 5 content, with M4 updated to 1)
$ws.Range("A4").Value = "Mac B. McGraw, Lindsay N. Kohler, Gabriel Q. Shaibi, Lawrence J. Mandarino, Dawn K. Coletta"
$ws.Range("B4").Value = "The University of Arizona College of Medicine, United States; Exos, United States; The University of Arizona, United States; Arizona State University, United States; The University of Arizona, United States; The University of Arizona, United States"
$ws.Range("C4").Value = "https://openalex.org/W4302286478"
$ws.Range("D4").Value = "A performance review of novel adiposity indices for assessing insulin resistance in a pediatric Latino population"
$ws.Range("E4").Value = "2022-10-06"
$ws.Range("F4").Value = "Frontiers in Pediatrics"
$ws.Range("G4").Value = "Frontiers Media"
$ws.Range("H4").Value = "https://doi.org/10.3389/fped.2022.1020901"
$ws.Range("I4").Value = "cc-by"
$ws.Range("J4").Value = "publishedVersion"
$ws.Range("K4").Value = "gold"
$ws.Range("M4").Value = "1"
$ws.Range("O4").Value = "https://pubmed.ncbi.nlm.nih.gov/36275055"
$ws.Range("P4").Value = "https://doi.org/10.3389/fped.2022.1020901"

# New row 5 (previously row 4 content, M5 stays 0)
$ws.Range("A5").Value = "Neusha Barakati, Rocio Zapata Bustos, Dawn K. Coletta, Paul Langlais, Lindsay N. Kohler, Moulun Luo, Janet L. Funk, Wayne T. Willis, Lawrence J. Mandarino"
$ws.Range("B5").Value = "Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Center for Disparities in Diabetes, Obesity, and Metabolism, University of Arizona, Health Sciences, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona"
$ws.Range("C5").Value = "https://openalex.org/W4281290394"
$ws.Range("D5").Value = "Acetylation of Adenine Nucleotide Translocase, Fuel Selection, and Metabolic Flexibility in Human Skeletal Muscle"
$ws.Range("E5").Value = "2022-05-09"
$ws.Range("F5").Value = "medRxiv (Cold Spring Harbor Laboratory)"
$ws.Range("G5").Value = "Cold Spring Harbor Laboratory"
$ws.Range("H5").Value = "https://doi.org/10.1101/2022.05.05.22274505"
$ws.Range("I5").Value = "N/A"
$ws.Range("J5").Value = "submittedVersion"
$ws.Range("K5").Value = "green"
$ws.Range("M5").Value = "0"
$ws.Range("O5").Value = "NA"
$ws.Range("P5").Value = "https://doi.org/10.1101/2022.05.05.22274505"
